# Auto-generated edit script: apply updated market-price / profit values
# to the Behemoth_Profits workbook, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2501.275
$ws.Cells.Item(132, 9).Value = 2038.4286
$ws.Cells.Item(132, 11).Value = 6115.2858
$ws.Cells.Item(132, 13).Value = -3585.2858
$ws.Cells.Item(137, 8).Value = 7220.1665
$ws.Cells.Item(137, 9).Value = 1842.1428
$ws.Cells.Item(137, 10).Value = 14749.4
$ws.Cells.Item(137, 11).Value = 5526.428400000001
$ws.Cells.Item(137, 12).Value = 44248.2
$ws.Cells.Item(137, 13).Value = -2976.428400000001
$ws.Cells.Item(137, 14).Value = -49348.2
$ws.Cells.Item(138, 8).Value = 1995.1351
$ws.Cells.Item(138, 9).Value = 1051.2
$ws.Cells.Item(138, 10).Value = 3105.647
$ws.Cells.Item(138, 11).Value = 3153.6
$ws.Cells.Item(138, 12).Value = 9316.940999999999
$ws.Cells.Item(138, 13).Value = 1986.4
$ws.Cells.Item(138, 14).Value = -19596.941
$ws.Cells.Item(141, 8).Value = 3400
$ws.Cells.Item(141, 9).Value = 3000
$ws.Cells.Item(141, 10).Value = 3666.6667
$ws.Cells.Item(141, 11).Value = 9000
$ws.Cells.Item(141, 12).Value = 11000.0001
$ws.Cells.Item(141, 13).Value = -3820
$ws.Cells.Item(141, 14).Value = -21360.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 1000.3333
$ws.Cells.Item(12, 9).Value = 1000.3333
$ws.Cells.Item(12, 11).Value = 1000.3333
$ws.Cells.Item(12, 13).Value = -827.3333
$ws.Cells.Item(45, 8).Value = 2572.2778
$ws.Cells.Item(45, 9).Value = 2272.9092
$ws.Cells.Item(45, 11).Value = 2272.9092
$ws.Cells.Item(45, 13).Value = -1895.9092
$ws.Cells.Item(61, 8).Value = 15661750
$ws.Cells.Item(61, 9).Value = 22731776
$ws.Cells.Item(61, 10).Value = 107694
$ws.Cells.Item(61, 11).Value = 22731776
$ws.Cells.Item(61, 12).Value = 107694
$ws.Cells.Item(61, 13).Value = -22731564
$ws.Cells.Item(61, 14).Value = -108118
$ws.Cells.Item(74, 8).Value = 19246460
$ws.Cells.Item(74, 9).Value = 41667416
$ws.Cells.Item(74, 11).Value = 41667416
$ws.Cells.Item(74, 13).Value = -41666542
$ws.Cells.Item(77, 8).Value = 19246460
$ws.Cells.Item(77, 9).Value = 41667416
$ws.Cells.Item(77, 11).Value = 208337080
$ws.Cells.Item(77, 13).Value = -208332712
$ws.Cells.Item(97, 8).Value = 962.4666999999999
$ws.Cells.Item(97, 9).Value = 752.8461
$ws.Cells.Item(97, 11).Value = 752.8461
$ws.Cells.Item(97, 13).Value = -256.8461
$ws.Cells.Item(122, 8).Value = 1404.0625
$ws.Cells.Item(122, 9).Value = 1342.25
$ws.Cells.Item(122, 10).Value = 1589.5
$ws.Cells.Item(122, 11).Value = 4026.75
$ws.Cells.Item(122, 12).Value = 4768.5
$ws.Cells.Item(122, 13).Value = -1576.75
$ws.Cells.Item(122, 14).Value = -9668.5
$ws.Cells.Item(134, 8).Value = 78998
$ws.Cells.Item(134, 10).Value = 78998
$ws.Cells.Item(134, 12).Value = 78998
$ws.Cells.Item(134, 14).Value = -89138
$ws.Cells.Item(136, 8).Value = 15661750
$ws.Cells.Item(136, 9).Value = 22731776
$ws.Cells.Item(136, 10).Value = 107694
$ws.Cells.Item(136, 11).Value = 68195328
$ws.Cells.Item(136, 12).Value = 323082
$ws.Cells.Item(136, 13).Value = -68192778
$ws.Cells.Item(136, 14).Value = -328182

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 45537.707
$ws.Cells.Item(134, 9).Value = 2070.6
$ws.Cells.Item(134, 10).Value = 76585.64
$ws.Cells.Item(134, 11).Value = 6211.799999999999
$ws.Cells.Item(134, 12).Value = 229756.92
$ws.Cells.Item(134, 13).Value = -3676.799999999999
$ws.Cells.Item(134, 14).Value = -234826.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1347.091
$ws.Cells.Item(58, 9).Value = 811
$ws.Cells.Item(58, 10).Value = 1990.4
$ws.Cells.Item(58, 11).Value = 811
$ws.Cells.Item(58, 12).Value = 1990.4
$ws.Cells.Item(58, 13).Value = -608
$ws.Cells.Item(58, 14).Value = -2396.4
$ws.Cells.Item(132, 8).Value = 3278.818
$ws.Cells.Item(132, 9).Value = 3219.6667
$ws.Cells.Item(132, 10).Value = 3545
$ws.Cells.Item(132, 11).Value = 9659.000100000001
$ws.Cells.Item(132, 12).Value = 10635
$ws.Cells.Item(132, 13).Value = -7129.000100000001
$ws.Cells.Item(132, 14).Value = -15695
$ws.Cells.Item(136, 8).Value = 1347.091
$ws.Cells.Item(136, 9).Value = 811
$ws.Cells.Item(136, 10).Value = 1990.4
$ws.Cells.Item(136, 11).Value = 2433
$ws.Cells.Item(136, 12).Value = 5971.200000000001
$ws.Cells.Item(136, 13).Value = 117
$ws.Cells.Item(136, 14).Value = -11071.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1526.4546
$ws.Cells.Item(5, 9).Value = 865.8889
$ws.Cells.Item(5, 10).Value = 4499
$ws.Cells.Item(5, 11).Value = 2597.6667
$ws.Cells.Item(5, 12).Value = 13497
$ws.Cells.Item(5, 13).Value = -2485.6667
$ws.Cells.Item(5, 14).Value = -13721
$ws.Cells.Item(11, 8).Value = 2451.238
$ws.Cells.Item(11, 9).Value = 2513.8
$ws.Cells.Item(11, 11).Value = 7541.400000000001
$ws.Cells.Item(11, 13).Value = -7401.400000000001
$ws.Cells.Item(107, 8).Value = 444.70587
$ws.Cells.Item(107, 10).Value = 430.66666
$ws.Cells.Item(107, 12).Value = 1291.99998
$ws.Cells.Item(107, 14).Value = -5131.999980000001
$ws.Cells.Item(109, 8).Value = 766.1429000000001
$ws.Cells.Item(109, 9).Value = 559.2308
$ws.Cells.Item(109, 11).Value = 1677.6924
$ws.Cells.Item(109, 13).Value = -637.6924000000001
$ws.Cells.Item(112, 8).Value = 9216.5
$ws.Cells.Item(112, 9).Value = 6433.3335
$ws.Cells.Item(112, 11).Value = 19300.0005
$ws.Cells.Item(112, 13).Value = -18192.0005
$ws.Cells.Item(132, 8).Value = 1720.4286
$ws.Cells.Item(132, 9).Value = 1590.8572
$ws.Cells.Item(132, 11).Value = 14317.7148
$ws.Cells.Item(132, 13).Value = -11787.7148
$ws.Cells.Item(135, 8).Value = 1526.4546
$ws.Cells.Item(135, 9).Value = 865.8889
$ws.Cells.Item(135, 10).Value = 4499
$ws.Cells.Item(135, 11).Value = 7793.0001
$ws.Cells.Item(135, 12).Value = 40491
$ws.Cells.Item(135, 13).Value = -5258.0001
$ws.Cells.Item(135, 14).Value = -45561

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 4597.4
$ws.Cells.Item(97, 9).Value = 4333
$ws.Cells.Item(97, 10).Value = 4994
$ws.Cells.Item(97, 11).Value = 4333
$ws.Cells.Item(97, 12).Value = 4994
$ws.Cells.Item(97, 13).Value = -3837
$ws.Cells.Item(97, 14).Value = -5986
$ws.Cells.Item(102, 8).Value = 6294.5625
$ws.Cells.Item(102, 9).Value = 4691.826
$ws.Cells.Item(102, 11).Value = 4691.826
$ws.Cells.Item(102, 13).Value = -3069.826
$ws.Cells.Item(122, 8).Value = 1975.9445
$ws.Cells.Item(122, 9).Value = 1975.1177
$ws.Cells.Item(122, 11).Value = 5925.3531
$ws.Cells.Item(122, 13).Value = -3475.3531
$ws.Cells.Item(128, 8).Value = 82388.89999999999
$ws.Cells.Item(128, 10).Value = 82388.89999999999
$ws.Cells.Item(128, 12).Value = 82388.89999999999
$ws.Cells.Item(128, 14).Value = -92348.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 44379.04
$ws.Cells.Item(7, 9).Value = 3891.8333
$ws.Cells.Item(7, 11).Value = 3891.8333
$ws.Cells.Item(7, 13).Value = -3779.8333
$ws.Cells.Item(22, 8).Value = 2682.5
$ws.Cells.Item(22, 9).Value = 2698.8235
$ws.Cells.Item(22, 11).Value = 2698.8235
$ws.Cells.Item(22, 13).Value = -2403.8235
$ws.Cells.Item(27, 8).Value = 2682.5
$ws.Cells.Item(27, 9).Value = 2698.8235
$ws.Cells.Item(27, 11).Value = 2698.8235
$ws.Cells.Item(27, 13).Value = -2591.8235
$ws.Cells.Item(68, 8).Value = 1664.6666
$ws.Cells.Item(68, 9).Value = 1664.6666
$ws.Cells.Item(68, 11).Value = 1664.6666
$ws.Cells.Item(68, 13).Value = -915.6666
$ws.Cells.Item(71, 8).Value = 1664.6666
$ws.Cells.Item(71, 9).Value = 1664.6666
$ws.Cells.Item(71, 11).Value = 8323.333000000001
$ws.Cells.Item(71, 13).Value = -4579.333000000001
$ws.Cells.Item(100, 8).Value = 2501.5
$ws.Cells.Item(100, 9).Value = 2501.5
$ws.Cells.Item(100, 11).Value = 2501.5
$ws.Cells.Item(100, 13).Value = -1960.5
$ws.Cells.Item(122, 8).Value = 5447.92
$ws.Cells.Item(122, 9).Value = 4511.9414
$ws.Cells.Item(122, 10).Value = 7436.875
$ws.Cells.Item(122, 11).Value = 13535.8242
$ws.Cells.Item(122, 12).Value = 22310.625
$ws.Cells.Item(122, 13).Value = -11085.8242
$ws.Cells.Item(122, 14).Value = -27210.625
$ws.Cells.Item(126, 8).Value = 44379.04
$ws.Cells.Item(126, 9).Value = 3891.8333
$ws.Cells.Item(126, 11).Value = 11675.4999
$ws.Cells.Item(126, 13).Value = -9205.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 75000
$ws.Cells.Item(70, 9).Value = 36000
$ws.Cells.Item(70, 10).Value = 114000
$ws.Cells.Item(70, 11).Value = 36000
$ws.Cells.Item(70, 12).Value = 114000
$ws.Cells.Item(70, 13).Value = -35685
$ws.Cells.Item(70, 14).Value = -114630
$ws.Cells.Item(73, 8).Value = 75000
$ws.Cells.Item(73, 9).Value = 36000
$ws.Cells.Item(73, 10).Value = 114000
$ws.Cells.Item(73, 11).Value = 36000
$ws.Cells.Item(73, 12).Value = 114000
$ws.Cells.Item(73, 13).Value = -34908
$ws.Cells.Item(73, 14).Value = -116184
$ws.Cells.Item(92, 8).Value = 42000
$ws.Cells.Item(92, 10).Value = 42000
$ws.Cells.Item(92, 12).Value = 42000
$ws.Cells.Item(92, 14).Value = -46992
$ws.Cells.Item(122, 8).Value = 4991.0625
$ws.Cells.Item(122, 9).Value = 2286.75
$ws.Cells.Item(122, 10).Value = 7695.375
$ws.Cells.Item(122, 11).Value = 6860.25
$ws.Cells.Item(122, 12).Value = 23086.125
$ws.Cells.Item(122, 13).Value = -4410.25
$ws.Cells.Item(122, 14).Value = -27986.125
$ws.Cells.Item(132, 8).Value = 17087.924
$ws.Cells.Item(132, 9).Value = 2143
$ws.Cells.Item(132, 11).Value = 6429
$ws.Cells.Item(132, 13).Value = -3899

Write-Host "Applied 215 cell updates across 8 sheets."
